# The edit re-shuffles the per-row data (Fecha, Variedad, Calidad, Volumen,
# Precio minimo/maximo/promedio, Unidad de comercializacion, Origen,
# Precio $/Kg, Kg / unidad) across the existing data rows (2..15) while
# leaving the rest of each row (Mercado, Region, Producto, Categoria, etc.)
# untouched. This corresponds to a fixed permutation of source rows per
# destination row, derived from the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as a group, by column index.
$cols = @(4, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)

# destinationRow -> sourceRow (values read from sourceRow BEFORE this edit,
# written into destinationRow).
$mapping = @{
    2  = 5
    3  = 6
    4  = 8
    5  = 15
    6  = 14
    7  = 11
    8  = 12
    9  = 4
    10 = 7
    11 = 13
    12 = 2
    13 = 3
    14 = 9
    15 = 10
}

# Snapshot the original values for every involved cell before writing
# anything, since the mapping is a permutation (destinations overlap with
# sources).
$snapshot = @{}
foreach ($row in 2..15) {
    foreach ($col in $cols) {
        $snapshot["$row,$col"] = $ws.Cells.Item($row, $col).Value2
    }
}

foreach ($destRow in 2..15) {
    $srcRow = $mapping[$destRow]
    foreach ($col in $cols) {
        $ws.Cells.Item($destRow, $col).Value = $snapshot["$srcRow,$col"]
    }
}
